# Added modifications to Distributed Energy 2030 scenario.
# Capacity, demand, fuelprices, EV fleet.
#
# Appends two new rows to the "constants" sheet for the "Distributed Energy"
# scenario, year 2030: one row with the C_0 coefficients, one with the C_1
# coefficients (all zero), following the same layout as the existing
# "National Trends"/2025 and "Distributed Energy"/2040 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Activate()

$data = New-Object 'object[,]' 2,27

# Row 6: Distributed Energy, 2030, C_0
$data[0,0]  = "Distributed Energy"
$data[0,1]  = 2030
$data[0,2]  = "C_0"
$data[0,3]  = 0      # AT
$data[0,4]  = 100    # BE
$data[0,5]  = 0      # CH
$data[0,6]  = 800    # DE
$data[0,7]  = 150    # DK_1
$data[0,8]  = 100    # DK_2
$data[0,9]  = 60     # EE
$data[0,10] = 0      # ES
$data[0,11] = 700    # FI
$data[0,12] = 2500   # FR
$data[0,13] = 1000   # GB
$data[0,14] = 30     # LT
$data[0,15] = 30     # LV
$data[0,16] = 600    # NL
$data[0,17] = 280    # NO_1
$data[0,18] = 200    # NO_2
$data[0,19] = 160    # NO_3
$data[0,20] = 120    # NO_4
$data[0,21] = 100    # NO_5
$data[0,22] = 1000   # PL
$data[0,23] = 200    # SE_1
$data[0,24] = 300    # SE_2
$data[0,25] = 500    # SE_3
$data[0,26] = 100    # SE_4

# Row 7: Distributed Energy, 2030, C_1
$data[1,0]  = "Distributed Energy"
$data[1,1]  = 2030
$data[1,2]  = "C_1"
for ($col = 3; $col -lt 27; $col++) {
    $data[1,$col] = 0
}

$ws.Range("A6:AA7").Value = $data

# Match the formatting used by the other data rows.
$ws.Range("A5:AA5").Copy()
$ws.Range("A6:AA7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M6").Select() | Out-Null
